$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A10").Value = "SendOnlySoftwareSerial"
$ws.Range("B10").Value = "Arduino Library"

$ws.Range("A11").Value = "PS2X_lib"
$ws.Range("B11").Value = "Arduino Library"

$ws.Range("A11").Select()
